# A new weekly price record for "Acelga" (Macroferia Regional de Talca) was
# published. It gets inserted at the top of the historical data block
# (row 154), pushing all the existing records down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 154; everything below (154:295) shifts
# down to (155:296), which also grows the sheet's dimension to A1:R296.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A154").Value = 5
$ws.Range("B154").Value = "Macroferia Regional de Talca"
$ws.Range("C154").Value = "Maule"
$ws.Range("D154").Value = 44790
$ws.Range("E154").Value = 7
$ws.Range("F154").Value = 100112009
$ws.Range("G154").Value = "Acelga"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 400
$ws.Range("K154").Value = 3500
$ws.Range("L154").Value = 3500
$ws.Range("M154").Value = 3500
$ws.Range("N154").Value = "$/docena de atados (4 kilos)"
$ws.Range("O154").Value = "Región del Maule"
$ws.Range("P154").Value = 875
$ws.Range("Q154").Value = 4
$ws.Range("R154").Value = "Hortaliza"
